$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial 45177 (2023-09-08) to 45178 (2023-09-09).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$ws.Range("C2:C$lastRow").Value = 45178
